# PlayerPerformance_4633.xlsx edit:
#  - insert a new "Player Info" sheet at the front
#  - rename MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" / "ODI Bowling"
#    and replace the scraped URL values with the bare match code
#  - drop the stray empty B22 cell on "ODI Batting"
#  - append a new "ODI Batting Extra" sheet at the end
#
# NOTE: worksheet object handles returned by Worksheets.Item(...) are
# position-bound in this host, not identity-bound - once a sheet is
# inserted/removed and the tab order shifts, a previously-fetched handle
# silently starts referring to whatever now sits at that old position. So
# every sheet handle below is (re)fetched by name immediately before it is
# used, and never reused across an intervening Worksheets.Add() call.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Player Info" sheet, inserted before the current first sheet
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($wb.Worksheets.Item("ODI Batting"))
$playerInfo.Name = "Player Info"

$playerInfoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $playerInfoHeaders.Length; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c)
    $cell.Value = $playerInfoHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

$playerInfo.Cells.Item(2, 1).NumberFormat = "@"
$playerInfo.Cells.Item(2, 1).Value = "4633"
$playerInfo.Cells.Item(2, 2).Value = "Najmul Hossain Shanto"
$playerInfo.Cells.Item(2, 3).Value = "Left Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Off Break"

# ---------------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE (header + values), drop B22
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")

$battingSheet.Cells.Item(1, 4).Value = "MATCH_CODE"

$battingCodes = @("4198", "4199", "4202", "4416", "4418", "4443", "4445", "4447", "4606", "4611", "4616", "4627", "4628", "4679", "4682", "4711", "4713", "4717", "4726", "4729", "4734")

$battingSheet.Range("D2:D22").NumberFormat = "@"
for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $battingSheet.Cells.Item($i + 2, 4).Value = $battingCodes[$i]
}

$battingSheet.Range("B22").ClearContents()

# ---------------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE (header + values)
# ---------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

$bowlingSheet.Cells.Item(1, 2).Value = "MATCH_CODE"

$bowlingCodes = @("4447", "4711")

$bowlingSheet.Range("B2:B3").NumberFormat = "@"
for ($i = 0; $i -lt $bowlingCodes.Length; $i++) {
    $bowlingSheet.Cells.Item($i + 2, 2).Value = $bowlingCodes[$i]
}

# ---------------------------------------------------------------------------
# 4. "ODI Batting Extra", appended after the last sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $extraHeaders.Length; $c++) {
    $cell = $extra.Cells.Item(1, $c)
    $cell.Value = $extraHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraData = @(
    ,@("4199", "2", "0", "0", "4.05%", "NO")
    ,@("4202", "2", "0", "0", "2.41%", "NO")
    ,@("4416", "3", "1", "2", "9.03%", "NO")
    ,@("4418", "", "", "", "", "NO")
    ,@("4443", "", "", "", "", "NO")
    ,@("4445", "3", "2", "0", "11.41%", "NO")
    ,@("4447", "3", "3", "0", "6.73%", "NO")
    ,@("4606", "", "", "", "", "NO")
    ,@("4611", "2", "2", "0", "17.86%", "NO")
    ,@("4616", "2", "0", "0", "0.56%", "NO")
    ,@("4627", "3", "5", "0", "13.10%", "NO")
    ,@("4628", "3", "0", "0", "", "NO")
    ,@("4679", "1", "0", "0", "", "NO")
    ,@("4682", "3", "3", "0", "7.75%", "NO")
    ,@("4711", "3", "6", "0", "27.75%", "NO")
    ,@("4713", "3", "0", "0", "", "NO")
    ,@("4717", "3", "5", "0", "21.54%", "NO")
    ,@("4726", "", "", "", "", "NO")
    ,@("4729", "", "", "", "", "NO")
    ,@("4734", "3", "", "", "", "NO")
)

# Columns A, C, D, E, F hold scraped text (even when it looks numeric, e.g.
# "0" or "4.05%"); only BATTING_POSITION (B) is a genuine number.
$extra.Range("A2:A21").NumberFormat = "@"
$extra.Range("C2:F21").NumberFormat = "@"

for ($r = 0; $r -lt $extraData.Length; $r++) {
    $row = $extraData[$r]
    $rowNum = $r + 2

    $extra.Cells.Item($rowNum, 1).Value = $row[0]

    if ($row[1] -ne "") {
        $extra.Cells.Item($rowNum, 2).Value = [double]$row[1]
    }

    if ($row[2] -ne "") {
        $extra.Cells.Item($rowNum, 3).Value = $row[2]
    }

    if ($row[3] -ne "") {
        $extra.Cells.Item($rowNum, 4).Value = $row[3]
    }

    if ($row[4] -ne "") {
        $extra.Cells.Item($rowNum, 5).Value = $row[4]
    }

    $extra.Cells.Item($rowNum, 6).Value = $row[5]
}

Write-Output "edit complete"
